$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force D2:E51 to Text format so numeric-looking strings (with trailing zeros,
# thousand-dot-like patterns, etc.) are preserved exactly as text, matching the
# original inline-string cell type used throughout this column range.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "30.047.55"
$ws.Range("E2").Value = "  -2.23%  "
$ws.Range("D3").Value = "2.106.70"
$ws.Range("E3").Value = "  -0.88%  "
$ws.Range("E4").Value = "  -1.05%  "
$ws.Range("D5").Value = "346.95"
$ws.Range("E5").Value = "  +2.25%  "
$ws.Range("E6").Value = "  -0.94%  "
$ws.Range("D7").Value = "0.5174"
$ws.Range("E7").Value = "  -1.99%  "
$ws.Range("D8").Value = "0.4435"
$ws.Range("E8").Value = "  -2.96%  "
$ws.Range("D9").Value = "0.09422"
$ws.Range("E9").Value = "  +3.39%  "
$ws.Range("D10").Value = "52.32"
$ws.Range("E10").Value = "  -3.22%  "
$ws.Range("D11").Value = "1.177"
$ws.Range("E11").Value = "  -0.04%  "
$ws.Range("D12").Value = "25.44"
$ws.Range("E12").Value = "  +3.70%  "
$ws.Range("D13").Value = "2.103.85"
$ws.Range("E13").Value = "  -1.13%  "
$ws.Range("D14").Value = "6.754"
$ws.Range("E14").Value = "  -1.52%  "
$ws.Range("D15").Value = "8.177"
$ws.Range("E15").Value = "  +0.64%  "
$ws.Range("D16").Value = "99.88"
$ws.Range("E16").Value = "  +2.38%  "
$ws.Range("D17").Value = "0.00001168"
$ws.Range("E17").Value = "  -0.47%  "
$ws.Range("E18").Value = "  -0.97%  "
$ws.Range("D19").Value = "20.80"
$ws.Range("E19").Value = "  +5.92%  "
$ws.Range("D20").Value = "0.06695"
$ws.Range("E20").Value = "  -0.47%  "
$ws.Range("E21").Value = "  -0.96%  "
$ws.Range("D22").Value = "6.250"
$ws.Range("E22").Value = "  -3.71%  "
$ws.Range("D23").Value = "30.125.54"
$ws.Range("D24").Value = "12.71"
$ws.Range("E24").Value = "  -2.80%  "
$ws.Range("D25").Value = "2.335"
$ws.Range("E25").Value = "  -2.13%  "
$ws.Range("D26").Value = "2.351.11"
$ws.Range("E26").Value = "  -1.01%  "
$ws.Range("D27").Value = "22.08"
$ws.Range("E27").Value = "  -2.14%  "
$ws.Range("D28").Value = "2.557"
$ws.Range("E28").Value = "  -0.16%  "
$ws.Range("D29").Value = "163.57"
$ws.Range("E29").Value = "  -1.31%  "
$ws.Range("D30").Value = "133.83"
$ws.Range("E30").Value = "  -0.98%  "
$ws.Range("D31").Value = "1.173"
$ws.Range("E31").Value = "  -3.27%  "
$ws.Range("D32").Value = "0.1061"
$ws.Range("E32").Value = "  -1.75%  "
$ws.Range("D33").Value = "1.644"
$ws.Range("E33").Value = "  -0.44%  "
$ws.Range("D34").Value = "6.256"
$ws.Range("E34").Value = "  -2.58%  "
$ws.Range("D35").Value = "3.955"
$ws.Range("E35").Value = "  -0.08%  "
$ws.Range("D36").Value = "6.237"
$ws.Range("E36").Value = "  +4.68%  "
$ws.Range("E37").Value = "  -4.17%  "
$ws.Range("D38").Value = "0.02570"
$ws.Range("E38").Value = "  -4.39%  "
$ws.Range("D39").Value = "0.06794"
$ws.Range("E39").Value = "  -1.56%  "
$ws.Range("D40").Value = "0.2294"
$ws.Range("E40").Value = "  -1.63%  "
$ws.Range("B41").Value = "TheSandbox"
$ws.Range("C41").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D41").Value = "0.6952"
$ws.Range("E41").Value = "  +0.28%  "
$ws.Range("B42").Value = "Aptos"
$ws.Range("C42").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D42").Value = "12.54"
$ws.Range("E42").Value = "  -1.05%  "
$ws.Range("D43").Value = "1.318"
$ws.Range("E43").Value = "  +3.82%  "
$ws.Range("D44").Value = "0.6646"
$ws.Range("E44").Value = "  +2.35%  "
$ws.Range("D45").Value = "14.32"
$ws.Range("E45").Value = "  -5.41%  "
$ws.Range("D46").Value = "2.300"
$ws.Range("E47").Value = "  -1.67%  "
$ws.Range("D48").Value = "0.00000000355"
$ws.Range("E48").Value = "  -3.85%  "
$ws.Range("D49").Value = "1.225"
$ws.Range("E49").Value = "  -3.03%  "
$ws.Range("D50").Value = "82.53"
$ws.Range("E50").Value = "  -1.24%  "
$ws.Range("D51").Value = "0.07215"
$ws.Range("E51").Value = "  -1.51%  "
